$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 89.5
$ws.Range("I5").Value = 89.5
$ws.Range("K5").Value = 89.5
$ws.Range("M5").Value = 25.5

$ws.Range("H41").Value = 487.63635
$ws.Range("I41").Value = 381.8
$ws.Range("J41").Value = 575.8333
$ws.Range("K41").Value = 381.8
$ws.Range("L41").Value = 575.8333
$ws.Range("M41").Value = 58.19999999999999
$ws.Range("N41").Value = -1455.8333

$ws.Range("H45").Value = 250
$ws.Range("J45").Value = 250
$ws.Range("L45").Value = 750
$ws.Range("N45").Value = -1134

$ws.Range("H62").Value = 29414846
$ws.Range("I62").Value = 32260880
$ws.Range("K62").Value = 32260880
$ws.Range("M62").Value = -32260256

$ws.Range("H65").Value = 29414846
$ws.Range("I65").Value = 32260880
$ws.Range("K65").Value = 161304400
$ws.Range("M65").Value = -161301280

$ws.Range("H98").Value = 22896.924
$ws.Range("I98").Value = 26799.182
$ws.Range("J98").Value = 1434.5
$ws.Range("K98").Value = 26799.182
$ws.Range("L98").Value = 1434.5
$ws.Range("M98").Value = -25301.182
$ws.Range("N98").Value = -4430.5

$ws.Range("H122").Value = 22896.924
$ws.Range("I122").Value = 26799.182
$ws.Range("J122").Value = 1434.5
$ws.Range("K122").Value = 80397.546
$ws.Range("L122").Value = 4303.5
$ws.Range("M122").Value = -77947.546
$ws.Range("N122").Value = -9203.5

$ws.Range("H129").Value = 55557320
$ws.Range("J129").Value = 142859600
$ws.Range("L129").Value = 428578800
$ws.Range("N129").Value = -428588800

$ws.Range("H132").Value = 2329452
$ws.Range("I132").Value = 4009.2896
$ws.Range("J132").Value = 20002816
$ws.Range("K132").Value = 12027.8688
$ws.Range("L132").Value = 60008448
$ws.Range("M132").Value = -9497.8688
$ws.Range("N132").Value = -60013508

$ws.Range("H137").Value = 6424.643
$ws.Range("I137").Value = 8812.611000000001
$ws.Range("J137").Value = 2126.3
$ws.Range("K137").Value = 26437.833
$ws.Range("L137").Value = 6378.900000000001
$ws.Range("M137").Value = -23887.833
$ws.Range("N137").Value = -11478.9

$ws.Range("H138").Value = 341243.03
$ws.Range("I138").Value = 1197917.1
$ws.Range("K138").Value = 3593751.3
$ws.Range("M138").Value = -3588611.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6373.171
$ws.Range("I32").Value = 6275.1484
$ws.Range("J32").Value = 10000
$ws.Range("K32").Value = 6275.1484
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = -5988.1484
$ws.Range("N32").Value = -10574

$ws.Range("H122").Value = 548484.3
$ws.Range("I122").Value = 2594.1702
$ws.Range("K122").Value = 7782.5106
$ws.Range("M122").Value = -5332.5106

$ws.Range("H132").Value = 5197.8
$ws.Range("I132").Value = 1738.75
$ws.Range("J132").Value = 6455.636
$ws.Range("K132").Value = 5216.25
$ws.Range("L132").Value = 19366.908
$ws.Range("M132").Value = -2686.25
$ws.Range("N132").Value = -24426.908

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 406.1111
$ws.Range("I22").Value = 431.875
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 431.875
$ws.Range("L22").Value = 200
$ws.Range("M22").Value = -81.875
$ws.Range("N22").Value = -900

$ws.Range("H31").Value = 12671
$ws.Range("J31").Value = 5973
$ws.Range("L31").Value = 5973
$ws.Range("N31").Value = -6563

$ws.Range("H34").Value = 12671
$ws.Range("J34").Value = 5973
$ws.Range("L34").Value = 5973
$ws.Range("N34").Value = -6377

$ws.Range("H64").Value = 50000
$ws.Range("J64").Value = 50000
$ws.Range("L64").Value = 50000
$ws.Range("N64").Value = -50496

$ws.Range("H67").Value = 50000
$ws.Range("J67").Value = 50000
$ws.Range("L67").Value = 50000
$ws.Range("N67").Value = -51716

$ws.Range("H107").Value = 55561930
$ws.Range("I107").Value = 83342750
$ws.Range("J107").Value = 288.66666
$ws.Range("K107").Value = 83342750
$ws.Range("L107").Value = 288.66666
$ws.Range("M107").Value = -83340830
$ws.Range("N107").Value = -4128.66666

$ws.Range("H132").Value = 1732.9524
$ws.Range("I132").Value = 1704.8422
$ws.Range("K132").Value = 5114.5266
$ws.Range("M132").Value = -2584.5266

$ws.Range("H134").Value = 6058.8667
$ws.Range("I134").Value = 4934.25
$ws.Range("K134").Value = 14802.75
$ws.Range("M134").Value = -12267.75

$ws.Range("H141").Value = 209720.89
$ws.Range("I141").Value = 82500
$ws.Range("J141").Value = 219898.56
$ws.Range("K141").Value = 82500
$ws.Range("L141").Value = 219898.56
$ws.Range("M141").Value = -77320
$ws.Range("N141").Value = -230258.56

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 13287.625
$ws.Range("J68").Value = 20730.2
$ws.Range("L68").Value = 62190.60000000001
$ws.Range("N68").Value = -63812.60000000001

$ws.Range("H71").Value = 13287.625
$ws.Range("J71").Value = 20730.2
$ws.Range("L71").Value = 186571.8
$ws.Range("N71").Value = -194683.8

$ws.Range("H107").Value = 1651.6666
$ws.Range("J107").Value = 1698.2858
$ws.Range("L107").Value = 5094.857400000001
$ws.Range("N107").Value = -8934.857400000001

$ws.Range("H117").Value = 397.6
$ws.Range("J117").Value = 503.33334
$ws.Range("L117").Value = 1510.00002
$ws.Range("N117").Value = -8394.000019999999

$ws.Range("H137").Value = 4366.5
$ws.Range("I137").Value = 1985.4348
$ws.Range("J137").Value = 8579.154
$ws.Range("K137").Value = 5956.3044
$ws.Range("L137").Value = 25737.462
$ws.Range("M137").Value = -856.3044
$ws.Range("N137").Value = -35937.462

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 15567.429
$ws.Range("I70").Value = 15997.5
$ws.Range("J70").Value = 15395.4
$ws.Range("K70").Value = 15997.5
$ws.Range("L70").Value = 15395.4
$ws.Range("M70").Value = -15727.5
$ws.Range("N70").Value = -15935.4

$ws.Range("H73").Value = 15567.429
$ws.Range("I73").Value = 15997.5
$ws.Range("J73").Value = 15395.4
$ws.Range("K73").Value = 15997.5
$ws.Range("L73").Value = 15395.4
$ws.Range("M73").Value = -15061.5
$ws.Range("N73").Value = -17267.4

$ws.Range("H122").Value = 10040.35
$ws.Range("I122").Value = 6755.1177
$ws.Range("K122").Value = 20265.3531
$ws.Range("M122").Value = -17815.3531

$ws.Range("H132").Value = 2348.1292
$ws.Range("I132").Value = 2360.7666
$ws.Range("K132").Value = 7082.2998
$ws.Range("M132").Value = -4552.2998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 1166.6666
$ws.Range("I21").Value = 225
$ws.Range("K21").Value = 225
$ws.Range("M21").Value = -51

$ws.Range("H40").Value = 18119.945
$ws.Range("I40").Value = 23719.908
$ws.Range("J40").Value = 9320
$ws.Range("K40").Value = 23719.908
$ws.Range("L40").Value = 9320
$ws.Range("M40").Value = -23583.908
$ws.Range("N40").Value = -9592

$ws.Range("H55").Value = 1034.8096
$ws.Range("I55").Value = 320.7857
$ws.Range("J55").Value = 2462.8572
$ws.Range("K55").Value = 320.7857
$ws.Range("L55").Value = 2462.8572
$ws.Range("M55").Value = -147.7857
$ws.Range("N55").Value = -2808.8572

$ws.Range("H122").Value = 5425.5186
$ws.Range("I122").Value = 4956.5654
$ws.Range("K122").Value = 14869.6962
$ws.Range("M122").Value = -12419.6962

$ws.Range("H132").Value = 515998.62
$ws.Range("I132").Value = 710999.5
$ws.Range("J132").Value = 4121.375
$ws.Range("K132").Value = 2132998.5
$ws.Range("L132").Value = 12364.125
$ws.Range("M132").Value = -2130468.5
$ws.Range("N132").Value = -17424.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 10710.695
$ws.Range("I132").Value = 12729.219
$ws.Range("K132").Value = 38187.657
$ws.Range("M132").Value = -35657.657
